$d = $word.ActiveDocument

# --- Edit 1: title paragraph "Employee Management System" becomes
#     "Payroll or Employee Management System", typed in as a new run
#     ("Payroll or ") ahead of the existing run, keeping the existing
#     run's formatting (bold, 36pt Times New Roman) ---
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleStart = $titleRange.Start

$insertPoint = $d.Range($titleStart, $titleStart)
$insertPoint.InsertBefore("Payroll or ")

# Force the newly typed text into its own run (rather than being
# silently merged back into the following run) by re-asserting its
# character formatting explicitly -- same Times New Roman / bold /
# 36pt formatting the rest of the title already uses.
$newRunRange = $d.Range($titleStart, $titleStart + 11)
$newRunRange.Font.Name = "Times New Roman"
$newRunRange.Font.NameFarEast = "Times New Roman"
$newRunRange.Font.NameBi = "Times New Roman"
$newRunRange.Font.Bold = $true
$newRunRange.Font.Size = 18

# --- Edit 2: "Create an employee's salary." collapses from three runs
#     (split as "Create an " / "employee's" / " salary.") into a single
#     run with the same text and formatting ---
$targetText = "Create an employee" + [char]0x2019 + "s salary."
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -eq $targetText) {
        $para.Range.Find.Execute($targetText, $false, $false, $false, $false, $false, $true, 1, $false, $targetText, 2)
        break
    }
}
